# Updates the cryptos price/volume table cell-by-cell to match the
# latest scrape (coinranking.com), per the Wed Mar 27 19:51:43 UTC 2024
# GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether the value must be
# forced to Text first (so Excel does not reinterpret a price like
# "569.72" or "5.00" or "0.0000300" as a floating point number and
# mangle its display / precision).
$updates = @(
    @{ Cell = 'D2'; Value = '68.711.79'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -1.69%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '3.493.48'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -2.10%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '569.72'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -1.52%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '182.59'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -3.35%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.614'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  -2.84%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '3.486.17'; ForceText = $false }
    @{ Cell = 'E8'; Value = '  -2.20%  '; ForceText = $false }
    @{ Cell = 'E9'; Value = '  +0.07%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.183'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +3.26%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.644'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -2.43%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '53.86'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -3.52%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '0.0000300'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  -0.95%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '9.41'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -1.96%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '4.043.77'; ForceText = $false }
    @{ Cell = 'E15'; Value = '  -2.46%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '19.21'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  -2.64%  '; ForceText = $false }
    @{ Cell = 'B17'; Value = 'WrappedBTC'; ForceText = $false }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; ForceText = $false }
    @{ Cell = 'D17'; Value = '68.622.82'; ForceText = $false }
    @{ Cell = 'E17'; Value = '  -1.61%  '; ForceText = $false }
    @{ Cell = 'B18'; Value = 'WrappedEther'; ForceText = $false }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; ForceText = $false }
    @{ Cell = 'D18'; Value = '3.474.08'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  -2.66%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '12.22'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -3.01%  '; ForceText = $false }
    @{ Cell = 'E20'; Value = '  -1.10%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '540.28'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +13.62%  '; ForceText = $false }
    @{ Cell = 'E22'; Value = '  -2.95%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '18.94'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -1.19%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '5.00'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -0.68%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '4.37'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -0.14%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '93.66'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -2.27%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '2.90'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -3.66%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '10.74'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -2.05%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '9.10'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -2.07%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '31.48'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -2.72%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '7.17'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -7.20%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '12.47'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  +2.41%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '64.59'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  -2.31%  '; ForceText = $false }
    @{ Cell = 'E34'; Value = '  -5.21%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '565.41'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -2.84%  '; ForceText = $false }
    @{ Cell = 'E36'; Value = '  -0.02%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '37.67'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -3.24%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.394'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  -0.12%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '2.97'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +4.42%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '0.0₃0763'; ForceText = $false }
    @{ Cell = 'E40'; Value = '  -3.99%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '3.10'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -3.93%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '3.32'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -3.95%  '; ForceText = $false }
    @{ Cell = 'E43'; Value = '  -4.04%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '3.252.01'; ForceText = $false }
    @{ Cell = 'E44'; Value = '  +0.82%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '3.48'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +2.91%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '2.96'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -3.70%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '0.0436'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -1.20%  '; ForceText = $false }
    @{ Cell = 'E48'; Value = '  -2.59%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '8.91'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -5.79%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.997'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -0.11%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '137.58'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  +1.98%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
